# Fixed a bug in Mask
# Reorders the data rows (rows 2-15 and rows 17-23) to reflect the corrected
# mask ordering. Rows 16, 24, 25 and the totals row 26 are unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(401, 9, 48, 67, 75, 45)
    3  = @(701, 3, 90, 45, 97, 15)
    4  = @(201, 9, 30, 15, 45, 30)
    5  = @(901, 16, 15, 45, 60, 60)
    6  = @(902, 1, 0, 0, 0, 0)
    7  = @(1203, 3, 15, 15, 15, 15)
    8  = @(101, 9, 30, 15, 60, 15)
    9  = @(1201, 2, 10, 10, 10, 10)
    10 = @(601, 9, 60, 67, 60, 42)
    11 = @(1202, 2, 10, 10, 10, 10)
    12 = @(501, 9, 52, 30, 75, 45)
    13 = @(801, 3, 67, 65, 52, 45)
    14 = @(1001, 18, 30, 75, 60, 72)
    15 = @(301, 6, 45, 30, 60, 45)
    17 = @(3, 0, 3, 3, 3, 3)
    18 = @(1, 0, 2, 2, 2, 2)
    19 = @(802, 0, 4, 5, 4, 0)
    20 = @(1101, 0, 15, 30, 30, 0)
    21 = @(2, 0, 2, 2, 2, 2)
    22 = @(602, 0, 0, 4, 0, 9)
    23 = @(402, 0, 0, 4, 0, 0)
}

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    $ws.Cells.Item($rowNum, 1).Value = $values[0]
    $ws.Cells.Item($rowNum, 2).Value = $values[1]
    $ws.Cells.Item($rowNum, 3).Value = $values[2]
    $ws.Cells.Item($rowNum, 4).Value = $values[3]
    $ws.Cells.Item($rowNum, 5).Value = $values[4]
    $ws.Cells.Item($rowNum, 6).Value = $values[5]
}
